$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Token (F) and EmailSent (G) columns entirely
$ws.Range("F1:G5").Delete()

# Add a new participant row who hasn't received the mail yet
$ws.Range("A6").Value = "tarun"
$ws.Range("B6").Value = "U654341"
$ws.Range("C6").Value = "2200039159@kluniversity.in"
$ws.Range("D6").Value = "Category-2"
$ws.Range("E6").Value = 2020
